# Update database and change read_price algorithm
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the period-header row (row 8) labels forward by one year,
#     dropping the oldest (1396/12) and introducing the newest (1401/12) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Shift the publish-date row (row 9) forward as well ---
$ws.Range("D9").Value = "1399-04-08 (8)"
$ws.Range("E9").Value = "1400-04-20 (8)"
$ws.Range("F9").Value = "1401-04-26 (9)"
$ws.Range("G9").Value = "1402-02-30 (8)"
$ws.Range("H9").Value = "1402-02-30"

# --- Zero-out (reset) the financial data body while the read_price
#     algorithm is reworked ---
$dataRows = 11,12,13,14,16,17,18,19,20,21,22,24,25,26,27
foreach ($r in $dataRows) {
    $ws.Range("D$r`:H$r").Value = 0
}

# Rows 15 and 23 use the "-" placeholder across the whole data range
$ws.Range("D15:H15").Value = "-"
$ws.Range("D23:H23").Value = "-"
